# Insert a new weekly price record as row 54 in the Mango price-history
# sheet. All existing rows from 54 downward shift down by one (54->55,
# ..., 139->140), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down by inserting a blank row at position 54.
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with the new record's values.
$ws.Range("A54").Value = 7
$ws.Range("B54").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C54").Value = "Ñuble"
$ws.Range("D54").Value = 45082
$ws.Range("E54").Value = 16
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100108
$ws.Range("H54").Value = "Tropicales y subtropicales"
$ws.Range("I54").Value = 100108002
$ws.Range("J54").Value = "Mango"
$ws.Range("K54").Value = "Sin especificar"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 60
$ws.Range("N54").Value = 9000
$ws.Range("O54").Value = 9000
$ws.Range("P54").Value = 9000
$ws.Range("Q54").Value = "$/bandeja 4 kilos"
$ws.Range("R54").Value = "Perú"
$ws.Range("S54").Value = 2250
$ws.Range("T54").Value = 4
